# Fixed update to excel issue
# Refresh the forecast for ASIN B0BF7FT26Z: the whole weekly series shifts
# one week forward (the oldest week drops off, a new week is appended at
# the end) and the forecast figures are recalculated. The Summary sheet is
# updated to match the refreshed Forecast Comparison data.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New weekly data for the "Forecast Comparison" sheet (rows 2-17).
# Columns: B=Week_Start_Date, D=MyForecast, E=Amazon Mean Forecast,
#          F=Amazon P70 Forecast, G=Amazon P80 Forecast, H=Amazon P90 Forecast
$forecastRows = @(
    @{ Row = 2;  Date = "2025-02-02"; D = 123; E = 99;  F = 115; G = 130; H = 151 },
    @{ Row = 3;  Date = "2025-02-09"; D = 104; E = 80;  F = 94;  G = 106; H = 126 },
    @{ Row = 4;  Date = "2025-02-16"; D = 110; E = 85;  F = 100; G = 113; H = 133 },
    @{ Row = 5;  Date = "2025-02-23"; D = 113; E = 87;  F = 103; G = 117; H = 140 },
    @{ Row = 6;  Date = "2025-03-02"; D = 116; E = 89;  F = 106; G = 122; H = 147 },
    @{ Row = 7;  Date = "2025-03-09"; D = 112; E = 86;  F = 102; G = 118; H = 142 },
    @{ Row = 8;  Date = "2025-03-16"; D = 116; E = 89;  F = 107; G = 125; H = 154 },
    @{ Row = 9;  Date = "2025-03-23"; D = 118; E = 91;  F = 111; G = 132; H = 166 },
    @{ Row = 10; Date = "2025-03-30"; D = 117; E = 90;  F = 107; G = 124; H = 152 },
    @{ Row = 11; Date = "2025-04-06"; D = 113; E = 87;  F = 105; G = 125; H = 156 },
    @{ Row = 12; Date = "2025-04-13"; D = 113; E = 87;  F = 105; G = 125; H = 155 },
    @{ Row = 13; Date = "2025-04-20"; D = 114; E = 88;  F = 107; G = 127; H = 159 },
    @{ Row = 14; Date = "2025-04-27"; D = 110; E = 85;  F = 103; G = 123; H = 154 },
    @{ Row = 15; Date = "2025-05-04"; D = 104; E = 80;  F = 97;  G = 117; H = 148 },
    @{ Row = 16; Date = "2025-05-11"; D = 107; E = 82;  F = 99;  G = 119; H = 151 },
    @{ Row = 17; Date = "2025-05-18"; D = 101; E = 78;  F = 95;  G = 114; H = 144 }
)

foreach ($entry in $forecastRows) {
    $r = $entry.Row
    # Leading apostrophe forces the date-like text to stay a plain text
    # string instead of being auto-converted into an Excel date value.
    $wsForecast.Cells.Item($r, 2).Value = "'" + $entry.Date
    $wsForecast.Cells.Item($r, 4).Value = $entry.D
    $wsForecast.Cells.Item($r, 5).Value = $entry.E
    $wsForecast.Cells.Item($r, 6).Value = $entry.F
    $wsForecast.Cells.Item($r, 7).Value = $entry.G
    $wsForecast.Cells.Item($r, 8).Value = $entry.H
}

# Updated metrics on the "Summary" sheet (column B), recomputed from the
# refreshed forecast above. Values that look like a plain number or a date
# get a leading apostrophe so they stay text cells, matching the workbook's
# original inline-string storage instead of becoming numeric/date cells.
$wsSummary.Range("B2").Value  = "2023-01-01 to 2025-01-26"
$wsSummary.Range("B3").Value  = "'0"
$wsSummary.Range("B4").Value  = "'302"
$wsSummary.Range("B6").Value  = "'162"
$wsSummary.Range("B7").Value  = "'74"
$wsSummary.Range("B8").Value  = "16592 units"
$wsSummary.Range("B9").Value  = "'1792"
$wsSummary.Range("B10").Value = "'912"
$wsSummary.Range("B11").Value = "'451"
$wsSummary.Range("B12").Value = "'123"
$wsSummary.Range("B14").Value = "'101"
$wsSummary.Range("B15").Value = "'2025-05-18"
